# "Done with remove duplicates"
# A leftover test row (a duplicate of row 8, "Abdukerim Ibrahim" / Uyghur /
# LO6DVTZLRK68528I / ... ) is appended as row 14. Unlike the original row 8,
# whose "version" column was pasted in as text, this one was typed in as a
# genuine number (5.9), which is why it lands as a numeric cell instead of
# a shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value2 = "Abdukerim Ibrahim"
$ws.Range("B14").Value2 = "Uyghur"
$ws.Range("C14").Value2 = "LO6DVTZLRK68528I"
$ws.Range("D14").Value2 = "Vivamus id faucibus velit, id posuere leo. Nunc aliquet sodales nunc a pulvinar. Nunc aliquet sodales nunc a pulvinar. Ut viverra quis eros eu tincidunt."
$ws.Range("E14").Value2 = 5.9

# Leave the freshly-typed row selected, the way it would be right after
# entering it on the keyboard.
$ws.Range("A14:E14").Select() | Out-Null
